$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Simple single-value cell replacements (row, new text)
$t.Cell(1, 1).Range.Text  = "0M"
$t.Cell(2, 1).Range.Text  = "0M"
$t.Cell(3, 1).Range.Text  = "0M"
$t.Cell(4, 1).Range.Text  = "133"
$t.Cell(6, 1).Range.Text  = "0.37486"
$t.Cell(7, 1).Range.Text  = "0.05817"
$t.Cell(8, 1).Range.Text  = "0.00924"
$t.Cell(9, 1).Range.Text  = "0.32385"
$t.Cell(10, 1).Range.Text = "0.37033"
$t.Cell(11, 1).Range.Text = "0.37438"
$t.Cell(12, 1).Range.Text = "2.00089"

# Collapse the multi-run tab-delimited cells down to a single value each
$t.Cell(44, 1).Range.Text = "92.63"
$t.Cell(45, 1).Range.Text = "2"
$t.Cell(46, 1).Range.Text = "27"
